$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_user")

# Update regcntr_id values per 2nd May data refresh
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Update the active view / selection to match saved workbook state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()
